# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps on the
# zh-cn and de-de report sheets, as part of regenerating the
# Handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2:E3").Value = "2016-03-30 11:02:18"
$wsZhCn.Range("H2:H3").Value = "2016-03-30 11:03:30"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2:E3").Value = "2016-03-30 11:02:32"
$wsDeDe.Range("H2:H3").Value = "2016-03-30 11:03:49"
